$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.321.81"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.91"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7098"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.95"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.0000"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07805"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3101"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.14"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.72%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08415"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.869.89"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.233"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7120"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.10%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.10"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008367"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +6.37%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.326.87"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.16%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.074"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "240.11"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.22"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.55%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.112.03"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9995"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.745"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1592"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.46"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.80%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.016"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.48"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.506"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.396"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.296"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.335"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05354"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.944"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.177"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7487"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.690"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01876"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.221.50"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +5.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.728"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.477"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8908"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "108.97"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +4.92%  "
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "72.34"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -1.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9997"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.015.10"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.03%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.796"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.60%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5197"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.42%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000122"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.427"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4321"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.50%  "
